$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 943-944, shifting the existing data (old rows 943-1023) down to 945-1025.
$ws.Rows("943:944").Insert(-4121)  # xlShiftDown

# Populate the two newly inserted rows with the new weekly price-report entries.
# row 943 (Conconina(o))
$ws.Cells.Item(943,1).Value2 = 11
$ws.Cells.Item(943,2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(943,3).Value2 = "Bíobío"
$ws.Cells.Item(943,4).Value2 = 45013
$ws.Cells.Item(943,5).Value2 = 8
$ws.Cells.Item(943,6).Value2 = 100112033
$ws.Cells.Item(943,7).Value2 = "Lechuga"
$ws.Cells.Item(943,8).Value2 = "Conconina(o)"
$ws.Cells.Item(943,9).Value2 = "Primera"
$ws.Cells.Item(943,10).Value2 = 220
$ws.Cells.Item(943,11).Value2 = 6000
$ws.Cells.Item(943,12).Value2 = 6500
$ws.Cells.Item(943,13).Value2 = 6273
$ws.Cells.Item(943,14).Value2 = "`$/caja 10 unidades"
$ws.Cells.Item(943,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(943,16).Value2 = 627
$ws.Cells.Item(943,17).Value2 = 10
$ws.Cells.Item(943,18).Value2 = "Hortaliza"

# row 944 (Escarola)
$ws.Cells.Item(944,1).Value2 = 11
$ws.Cells.Item(944,2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(944,3).Value2 = "Bíobío"
$ws.Cells.Item(944,4).Value2 = 45013
$ws.Cells.Item(944,5).Value2 = 8
$ws.Cells.Item(944,6).Value2 = 100112033
$ws.Cells.Item(944,7).Value2 = "Lechuga"
$ws.Cells.Item(944,8).Value2 = "Escarola"
$ws.Cells.Item(944,9).Value2 = "Primera"
$ws.Cells.Item(944,10).Value2 = 220
$ws.Cells.Item(944,11).Value2 = 6000
$ws.Cells.Item(944,12).Value2 = 6500
$ws.Cells.Item(944,13).Value2 = 6227
$ws.Cells.Item(944,14).Value2 = "`$/caja 15 unidades"
$ws.Cells.Item(944,15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(944,16).Value2 = 415
$ws.Cells.Item(944,17).Value2 = 15
$ws.Cells.Item(944,18).Value2 = "Hortaliza"
